$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.562.44'
$ws.Range('E2').Value = '  +5.46%  '
$ws.Range('D3').Value = '1.917.67'
$ws.Range('E3').Value = '  +3.85%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  -0.75%  '
$ws.Range('D5').Value = "'335.08"
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  +3.29%  '
$ws.Range('D8').Value = "'0.4128"
$ws.Range('E8').Value = '  +5.95%  '
$ws.Range('D9').Value = "'48.13"
$ws.Range('E9').Value = '  +1.59%  '
$ws.Range('D10').Value = "'0.08040"
$ws.Range('E10').Value = '  +3.83%  '
$ws.Range('D11').Value = "'1.014"
$ws.Range('E11').Value = '  +3.84%  '
$ws.Range('D12').Value = "'22.46"
$ws.Range('E12').Value = '  +6.31%  '
$ws.Range('D13').Value = '1.974.20'
$ws.Range('E13').Value = '  +6.21%  '
$ws.Range('D14').Value = "'6.012"
$ws.Range('E14').Value = '  +4.21%  '
$ws.Range('D15').Value = "'7.192"
$ws.Range('E15').Value = '  +3.12%  '
$ws.Range('D16').Value = "'89.99"
$ws.Range('E16').Value = '  +3.59%  '
$ws.Range('D17').Value = "'1.004"
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').Value = "'0.00001037"
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').Value = "'0.06598"
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').Value = "'17.86"
$ws.Range('E20').Value = '  +5.77%  '
$ws.Range('D21').Value = "'1.002"
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').Value = '29.553.63'
$ws.Range('E22').Value = '  +5.34%  '
$ws.Range('D23').Value = "'5.569"
$ws.Range('E23').Value = '  +5.09%  '
$ws.Range('D24').Value = "'11.64"
$ws.Range('E24').Value = '  +9.99%  '
$ws.Range('E25').Value = '  -2.58%  '
$ws.Range('D26').Value = '2.190.43'
$ws.Range('E26').Value = '  +5.51%  '
$ws.Range('D27').Value = "'156.93"
$ws.Range('E27').Value = '  +1.30%  '
$ws.Range('D28').Value = "'19.92"
$ws.Range('E28').Value = '  +4.30%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = "'2.145"
$ws.Range('E29').Value = '  +5.72%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = "'5.739"
$ws.Range('E30').Value = '  +10.01%  '
$ws.Range('D31').Value = "'117.56"
$ws.Range('E31').Value = '  +1.29%  '
$ws.Range('D32').Value = "'1.070"
$ws.Range('E32').Value = '  +15.12%  '
$ws.Range('D33').Value = "'0.09481"
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('D34').Value = "'1.434"
$ws.Range('E34').Value = '  +5.05%  '
$ws.Range('D35').Value = "'5.431"
$ws.Range('E35').Value = '  +5.03%  '
$ws.Range('D36').Value = "'3.529"
$ws.Range('E36').Value = '  -2.37%  '
$ws.Range('D37').Value = "'0.06143"
$ws.Range('E37').Value = '  +2.60%  '
$ws.Range('D38').Value = "'0.02270"
$ws.Range('E38').Value = '  +4.31%  '
$ws.Range('D39').Value = "'8.450"
$ws.Range('E39').Value = '  +4.06%  '
$ws.Range('D40').Value = "'1.181"
$ws.Range('E40').Value = '  +3.36%  '
$ws.Range('D41').Value = "'0.5902"
$ws.Range('E41').Value = '  +4.63%  '
$ws.Range('D42').Value = "'0.1846"
$ws.Range('E42').Value = '  +3.63%  '
$ws.Range('D43').Value = "'10.20"
$ws.Range('E43').Value = '  +3.31%  '
$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'1.261"
$ws.Range('E44').Value = '  +1.94%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'2.337"
$ws.Range('E45').Value = '  +3.44%  '
$ws.Range('D46').Value = "'0.07509"
$ws.Range('E46').Value = '  +4.79%  '
$ws.Range('D47').Value = "'0.5588"
$ws.Range('E47').Value = '  +4.39%  '
$ws.Range('D48').Value = "'12.21"
$ws.Range('E48').Value = '  +4.28%  '
$ws.Range('D49').Value = "'1.938"
$ws.Range('E49').Value = '  +4.26%  '
$ws.Range('D50').Value = "'113.18"
$ws.Range('E50').Value = '  +3.41%  '
$ws.Range('D51').Value = "'0.2983"
$ws.Range('E51').Value = '  +13.77%  '
